$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub6 = [char]0x2086

$updates = @(
    @{Row=2;  D="65.824.28";  E="  +0.16%  "},
    @{Row=3;  D="2.664.83";   E="  -0.49%  "},
    @{Row=4;  D=$null;        E="  +0.01%  "},
    @{Row=5;  D="598.55";     E="  -0.38%  "},
    @{Row=6;  D="158.35";     E="  +1.00%  "},
    @{Row=7;  D="0.654";      E="  +4.93%  "},
    @{Row=8;  D=$null;        E="  +0.01%  "},
    @{Row=9;  D=$null;        E="  -2.69%  "},
    @{Row=10; D=$null;        E="  +0.62%  "},
    @{Row=11; D="5.86";       E="  -0.46%  "},
    @{Row=12; D=$null;        E="  +1.53%  "},
    @{Row=13; D="29.03";      E="  -1.33%  "},
    @{Row=14; D=$null;        E="  -2.49%  "},
    @{Row=15; D="3.143.50";   E="  -0.50%  "},
    @{Row=16; D="65.687.09";  E="  +0.18%  "},
    @{Row=17; D="2.664.03";   E="  -0.26%  "},
    @{Row=18; D="12.61";      E="  -2.53%  "},
    @{Row=19; D="4.81";       E="  -0.06%  "},
    @{Row=20; D="7.50";       E="  -1.22%  "},
    @{Row=21; D="351.48";     E="  -0.33%  "},
    @{Row=22; D=$null;        E="  -0.08%  "},
    @{Row=23; D=$null;        E="  +0.13%  "},
    @{Row=24; D=$null;        E="  +11.75%  "},
    @{Row=25; D=$null;        E="  +0.54%  "},
    @{Row=26; D="9.65";       E="  -0.40%  "},
    @{Row=27; D=$null;        E="  +1.24%  "},
    @{Row=28; D="570.09";     E="  +7.36%  "},
    @{Row=29; D="8.21";       E="  +1.43%  "},
    @{Row=30; D=$null;        E="  -2.56%  "},
    @{Row=31; D=$null;        E="  -0.19%  "},
    @{Row=32; D="2.16";       E="  +0.62%  "},
    @{Row=33; D="1.82";       E="  +3.77%  "},
    @{Row=34; D="6.71";       E="  +3.43%  "},
    @{Row=35; D="5.57";       E="  +1.49%  "},
    @{Row=36; D=$null;        E="  -0.09%  "},
    @{Row=37; D="20.65";      E="  +0.28%  "},
    @{Row=38; D=$null;        E="  -0.01%  "},
    @{Row=39; D="1.96";       E="  +0.45%  "},
    @{Row=40; D="154.62";     E="  -2.42%  "},
    @{Row=41; D="161.84";     E="  -1.83%  "},
    @{Row=42; D="4.11";       E="  -0.87%  "},
    @{Row=43; D="0.0619";     E="  +1.30%  "},
    @{Row=44; D="2.32";       E="  -0.26%  "},
    @{Row=45; D="23.13";      E="  +0.94%  "},
    @{Row=46; D=$null;        E="  +0.36%  "},
    @{Row=47; D="0.0258";     E="  -0.27%  "},
    @{Row=48; D="0.103";      E="  +2.05%  "},
    @{Row=49; D="19.86";      E="  -1.21%  "},
    @{Row=50; D="0.0${sub6}0246"; E="  -5.36%  "},
    @{Row=51; D="0.817";      E="  -0.15%  "}
)

# All D/E cells in this sheet are text strings (coinranking.com "Price" /
# "Volume(1h)" columns, e.g. "65.824.28" or "  +0.16%  "), never real
# numbers. Excel, however, auto-converts a plainly-numeric-looking string
# (single decimal point, no thousands separators, e.g. "598.55") assigned
# via Value/Value2 into a genuine number. Pre-format those specific cells
# as Text so the written value round-trips as a string, matching the
# source data (prices with thousands separators like "65.824.28", or
# non-numeric text, are left alone since Excel already keeps those as text).
$forceTextRows = @(5,6,7,11,13,18,19,20,21,26,28,29,32,33,34,35,37,39,40,41,42,43,44,45,47,48,49,51)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $dcell = $ws.Cells.Item($r, 4)
        if ($forceTextRows -contains $r) {
            $dcell.NumberFormat = "@"
        }
        $dcell.Value2 = $u.D
    }
    $ws.Cells.Item($r, 5).Value2 = $u.E
}
